# Update vanwege 25e toets (dus nu 26 toetsmodules met menu-knop)
#
# The "Toewijzingen" pin-mapping sheet gets re-jiggled: the standalone
# "Menu" key becomes "Tmenu" (test-menu), the old SH/LDn1..3 trio collapses
# to a single shared "SH/LDn" line, and the freed-up rows get re-used for
# the new 25th test key ("T25") and a spare "NC" (not connected) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Toewijzingen")

$YELLOW = 65535     # RGB(255,255,0)  -> fillId2 "FFFFFF00"
$GREEN  = 5296274   # RGB(146,208,80) -> fillId4 "FF92D050"

# Row 5: "SH/LDn1" -> "SH/LDn" (fill/style unchanged - stays green)
$ws.Range("A5").Value2 = "SH/LDn"
$ws.Range("A5").Interior.Color = $GREEN

# Row 4: "Menu" -> "Tmenu" (fill/style unchanged - stays yellow)
$ws.Range("F4").Value2 = "Tmenu"
$ws.Range("F4").Interior.Color = $YELLOW

# Row 6: "SH/LDn2" -> "DIN", "Shift" -> "T25" (both stay yellow)
$ws.Range("A6").Value2 = "DIN"
$ws.Range("A6").Interior.Color = $YELLOW
$ws.Range("B6").Value2 = "T25"
$ws.Range("B6").Interior.Color = $YELLOW

# Row 6: F6/G6 text unchanged ("BatMon"/"ADC" stay as-is, still yellow)

# Row 7: "SH/LDn3" -> "NC", "Shift" -> "NC" (both switch yellow -> green)
$ws.Range("A7").Value2 = "NC"
$ws.Range("A7").Interior.Color = $GREEN
$ws.Range("B7").Value2 = "NC"
$ws.Range("B7").Interior.Color = $GREEN

# Row 8: text stays "CS" (underlying shared-string slot shifted only)
$ws.Range("A8").Value2 = "CS"
